$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Row 36: relabel existing "Sol" placeholder into the real solution name,
#     and fill in the rest of the row with the new data point. ---
$ws.Cells.Item(36, 1).Value = "Solution 33: BRAM internal arrays"
$ws.Cells.Item(36, 2).Value = 64
$ws.Cells.Item(36, 3).Value = 10817.62
$ws.Cells.Item(36, 4).Value = "Main, MULT_I II=8"
$ws.Cells.Item(36, 5).Value = "All the rest"
$ws.Cells.Item(36, 6).Value = "RTL"
$ws.Cells.Item(36, 7).Value = "Set transpose in intialization"
$ws.Cells.Item(36, 8).Formula = '=B36*C36'
$ws.Cells.Item(36, 9).Formula = '=H36/MIN($H$5:$H$110)'

# --- Row 37: second new solution row (same solution/category, different
#     measurement -- internal arrays mapped to registers). ---
$ws.Cells.Item(37, 1).Value = "Solution 33: BRAM internal arrays"
$ws.Cells.Item(37, 2).Value = 64
$ws.Cells.Item(37, 3).Value = 8839.5400000000009
$ws.Cells.Item(37, 4).Value = "Main, MULT_I II=8"
$ws.Cells.Item(37, 5).Value = "All the rest"
$ws.Cells.Item(37, 6).Value = "RTL"
$ws.Cells.Item(37, 7).Value = "Internal arrays mapped to registers"
$ws.Cells.Item(37, 8).Formula = '=B37*C37'
$ws.Cells.Item(37, 9).Formula = '=H37/MIN($H$5:$H$110)'

# --- Update the saved view state to match where the author left the cursor. ---
[void]$ws.Range("G38").Select()
$win = $excel.ActiveWindow
$win.ScrollRow = 34
$win.ScrollColumn = 1
